# The workbook tracks daily wholesale prices for "Berenjena" (eggplant).
# This commit adds one new day of data at the top of the date-ordered block
# (old row 228), pushing the existing rows 228-331 down by one (to 229-332).
#
# Net effect:
#   - Insert a new row at position 228, shifting rows 228..331 down to 229..332
#     (this naturally grows the used range from A1:R331 to A1:R332).
#   - Populate the new row 228 with the new day's record. Most descriptive
#     columns (market/region/category/etc.) are identical to the rest of the
#     block, only the date and the measured values differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 228; rows 228-331 shift to 229-332.
$ws.Rows.Item(228).Insert()

# Fill in the new row 228 with the new record.
$ws.Cells.Item(228, 1).Value  = 6
$ws.Cells.Item(228, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(228, 3).Value  = "Metropolitana"
$ws.Cells.Item(228, 4).Value  = 45119
$ws.Cells.Item(228, 5).Value  = 13
$ws.Cells.Item(228, 6).Value  = 100112001
$ws.Cells.Item(228, 7).Value  = "Berenjena"
$ws.Cells.Item(228, 8).Value  = "Sin especificar"
$ws.Cells.Item(228, 9).Value  = "Primera"
$ws.Cells.Item(228, 10).Value = 250
$ws.Cells.Item(228, 11).Value = 7000
$ws.Cells.Item(228, 12).Value = 8000
$ws.Cells.Item(228, 13).Value = 7600
$ws.Cells.Item(228, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(228, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(228, 16).Value = 152
$ws.Cells.Item(228, 17).Value = 50
$ws.Cells.Item(228, 18).Value = "Hortaliza"
